# Add a new per-user worksheet "workout_data_mahad123" to the workout-tracker
# workbook, mirroring the existing "workout_data_tvisha" sheet (same headers,
# same 20 rows of sample data, same column styling) -- this is how a new user
# sheet gets seeded in this workbook.

$wb = $excel.ActiveWorkbook

# Remember which sheet was active/selected so we can restore that after the
# copy (Excel activates the newly created sheet by default).
$originalActiveSheetName = $wb.ActiveSheet.Name

$template = $wb.Worksheets.Item("workout_data_tvisha")

# Copy the template sheet, placing the new copy immediately after it.
$template.Copy([System.Reflection.Missing]::Value, $template)

# The freshly copied sheet is now the last sheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "workout_data_mahad123"

# Restore the original active sheet/selection.
$wb.Worksheets.Item($originalActiveSheetName).Activate()

Write-Output ("Sheets now: " + (($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "))
